$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a string value to a cell while forcing it to stay a TEXT
# cell (matching the source file, where every Price/Volume cell is an
# inline string) even when the string looks like a plain number
# (e.g. "1.004"). Excel's normal type-inference would otherwise store
# such values as numeric cells. Temporarily flipping the number format
# to "@" (Text) before the assignment prevents that inference, and
# resetting the style back to "Normal" afterwards avoids leaving any
# stray formatting behind on the cell.
function Set-CellText($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-CellText $ws.Range("D2") "26.426.05"
$ws.Range("E2").Value = "  -3.80%  "

# Row 3 - Ethereum
Set-CellText $ws.Range("D3") "1.769.94"
$ws.Range("E3").Value = "  -3.03%  "

# Row 4 - TetherUSD
Set-CellText $ws.Range("D4") "1.004"

# Row 5 - USDC
$ws.Range("E5").Value = "  +0.09%  "

# Row 6 - BNB
Set-CellText $ws.Range("D6") "306.38"
$ws.Range("E6").Value = "  -2.10%  "

# Row 7 - XRP
Set-CellText $ws.Range("D7") "0.4303"
$ws.Range("E7").Value = "  +1.10%  "

# Row 8 - Cardano
Set-CellText $ws.Range("D8") "0.3661"
$ws.Range("E8").Value = "  +1.37%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.17%  "

# Row 10 - Polygon
Set-CellText $ws.Range("D10") "0.8487"
$ws.Range("E10").Value = "  -1.63%  "

# Row 11 - Solana
Set-CellText $ws.Range("D11") "20.37"
$ws.Range("E11").Value = "  -1.01%  "

# Row 12 - WrappedEther
Set-CellText $ws.Range("D12") "1.787.97"
$ws.Range("E12").Value = "  -8.27%  "

# Row 13 - Chainlink
Set-CellText $ws.Range("D13") "6.428"
$ws.Range("E13").Value = "  -0.74%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -2.76%  "

# Row 15 - TRON
Set-CellText $ws.Range("D15") "0.06946"
$ws.Range("E15").Value = "  +0.22%  "

# Row 16 - BinanceUSD
$ws.Range("E16").Value = "  -0.05%  "

# Row 17 - Litecoin
Set-CellText $ws.Range("D17") "79.29"
$ws.Range("E17").Value = "  -1.87%  "

# Row 18 - ShibaInu
Set-CellText $ws.Range("D18") "0.000008677"
$ws.Range("E18").Value = "  -2.78%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.15%  "

# Row 20 - Avalanche
Set-CellText $ws.Range("D20") "15.04"
$ws.Range("E20").Value = "  -2.40%  "

# Row 21 - WrappedBTC
Set-CellText $ws.Range("D21") "26.444.20"
$ws.Range("E21").Value = "  -5.53%  "

# Row 22 - Uniswap
Set-CellText $ws.Range("D22") "5.098"
$ws.Range("E22").Value = "  -0.53%  "

# Row 23 - Cosmos
Set-CellText $ws.Range("D23") "11.23"
$ws.Range("E23").Value = "  +3.26%  "

# Row 24 - WrappedliquidstakedEther2.0
Set-CellText $ws.Range("D24") "2.006.46"
$ws.Range("E24").Value = "  -7.51%  "

# Row 25 - Monero
Set-CellText $ws.Range("D25") "152.07"
$ws.Range("E25").Value = "  -2.15%  "

# Row 26 - Toncoin
Set-CellText $ws.Range("D26") "1.872"
$ws.Range("E26").Value = "  -6.00%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  -3.27%  "

# Row 28 - InternetComputer(DFINITY)
Set-CellText $ws.Range("D28") "5.087"
$ws.Range("E28").Value = "  -0.83%  "

# Row 29 - BitcoinCash
Set-CellText $ws.Range("D29") "114.47"
$ws.Range("E29").Value = "  +0.27%  "

# Row 30 - LidoDAOToken
Set-CellText $ws.Range("D30") "1.757"
$ws.Range("E30").Value = "  -1.95%  "

# Row 31 - Stellar
Set-CellText $ws.Range("D31") "0.08966"
$ws.Range("E31").Value = "  +0.78%  "

# Row 32 - ImmutableX
Set-CellText $ws.Range("D32") "0.7247"
$ws.Range("E32").Value = "  -2.99%  "

# Row 33 - ARBITRUM
Set-CellText $ws.Range("D33") "1.112"
$ws.Range("E33").Value = "  -0.59%  "

# Row 34 - Filecoin
Set-CellText $ws.Range("D34") "4.327"
$ws.Range("E34").Value = "  -4.74%  "

# Rows 35 & 36 swap places: HuobiToken <-> Frax, each keeping its own
# freshly-updated price/volume reading.
$ws.Range("B35").Value = "Frax"
$ws.Range("C35").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-CellText $ws.Range("D35") "1.003"
$ws.Range("E35").Value = "  +0.14%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-CellText $ws.Range("D36") "2.745"
$ws.Range("E36").Value = "  -8.05%  "

# Row 37 - TrustWalletToken
Set-CellText $ws.Range("D37") "1.078"
$ws.Range("E37").Value = "  -0.59%  "

# Row 38 - Hedera
Set-CellText $ws.Range("D38") "0.05151"
$ws.Range("E38").Value = "  -2.09%  "

# Row 39 - VeChain
Set-CellText $ws.Range("D39") "0.01889"
$ws.Range("E39").Value = "  -1.79%  "

# Row 40 - TheSandbox
Set-CellText $ws.Range("D40") "0.4923"
$ws.Range("E40").Value = "  -3.01%  "

# Row 41 - Algorand
$ws.Range("E41").Value = "  -3.04%  "

# Row 42 - MXToken
$ws.Range("E42").Value = "  -7.65%  "

# Row 43 - FraxShare
Set-CellText $ws.Range("D43") "6.265"
$ws.Range("E43").Value = "  -1.87%  "

# Row 44 - Aptos
Set-CellText $ws.Range("D44") "8.006"
$ws.Range("E44").Value = "  -4.07%  "

# Row 45 - Quant
Set-CellText $ws.Range("D45") "104.80"
$ws.Range("E45").Value = "  -1.51%  "

# Row 46 - PaxDollar
Set-CellText $ws.Range("D46") "1.002"
$ws.Range("E46").Value = "  +0.15%  "

# Row 47 - EnergySwap
Set-CellText $ws.Range("D47") "10.14"
$ws.Range("E47").Value = "  -2.81%  "

# Rows 48 & 49 swap places: Cronos <-> Decentraland, each keeping its own
# freshly-updated price/volume reading.
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-CellText $ws.Range("D48") "0.4493"
$ws.Range("E48").Value = "  -4.08%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-CellText $ws.Range("D49") "0.06194"
$ws.Range("E49").Value = "  -4.07%  "

# Row 50 - NEARProtocol
$ws.Range("E50").Value = "  -1.49%  "

# Row 51 - RenderToken
Set-CellText $ws.Range("D51") "1.737"
$ws.Range("E51").Value = "  +2.71%  "
